$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking target cells to stay text cells (matches the
# original inline-string cell type instead of being auto-coerced to numbers).
$textCells = @(
    'D5',
    'D6',
    'D14',
    'D20',
    'D21',
    'D24',
    'D25',
    'D26',
    'D30',
    'D32',
    'D33',
    'D34',
    'D35',
    'D37',
    'D38',
    'D39',
    'D41',
    'D42',
    'D45',
    'D47',
    'D48',
    'D49'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values.
$ws.Range('D2').Value = '61.304.93'
$ws.Range('E2').Value = '  +0.96%  '
$ws.Range('D3').Value = '2.374.72'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '551.42'
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('D6').Value = '139.67'
$ws.Range('E6').Value = '  +1.69%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '2.376.13'
$ws.Range('E10').Value = '  +3.60%  '
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('E12').Value = '  +2.22%  '
$ws.Range('E13').Value = '  +3.04%  '
$ws.Range('D14').Value = '25.48'
$ws.Range('E14').Value = '  +2.76%  '
$ws.Range('E15').Value = '  +5.18%  '
$ws.Range('D16').Value = '2.805.15'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '61.345.25'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('D18').Value = '2.379.87'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('E19').Value = '  +3.35%  '
$ws.Range('D20').Value = '4.15'
$ws.Range('E20').Value = '  +2.29%  '
$ws.Range('D21').Value = '320.56'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '1.75'
$ws.Range('E24').Value = '  -6.20%  '
$ws.Range('D25').Value = '64.33'
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('D26').Value = '8.84'
$ws.Range('E26').Value = '  +7.61%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '2.493.41'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E29').Value = '  +3.13%  '
$ws.Range('D30').Value = '517.43'
$ws.Range('E30').Value = '  +3.43%  '
$ws.Range('D31').Value = '0.0₃0899'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.39'
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = '0.149'
$ws.Range('E33').Value = '  +3.52%  '
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').Value = '  +3.16%  '
$ws.Range('D35').Value = '1.53'
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '5.51'
$ws.Range('E37').Value = '  +5.65%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '4.68'
$ws.Range('E38').Value = '  +2.92%  '
$ws.Range('D39').Value = '1.89'
$ws.Range('E39').Value = '  +5.25%  '
$ws.Range('D41').Value = '18.49'
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').Value = '146.87'
$ws.Range('E42').Value = '  +5.77%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('D45').Value = '147.79'
$ws.Range('E45').Value = '  +6.72%  '
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('D47').Value = '3.61'
$ws.Range('E47').Value = '  +2.19%  '
$ws.Range('D48').Value = '0.0524'
$ws.Range('E48').Value = '  +2.54%  '
$ws.Range('D49').Value = '19.62'
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('E51').Value = '  +1.23%  '
